$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the data set. It belongs
# chronologically between the existing row 22 (2020-12-04) and the old
# row 23 (2020-12-07), so insert a fresh row at position 23, pushing the
# old rows 23-29 down to 24-30.
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new record's data.
$ws.Cells.Item(23, 1).Value = 5
$ws.Cells.Item(23, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(23, 3).Value = "Maule"
$ws.Cells.Item(23, 4).Value = 44524
$ws.Cells.Item(23, 5).Value = 7
$ws.Cells.Item(23, 6).Value = "Fruta"
$ws.Cells.Item(23, 7).Value = 100101
$ws.Cells.Item(23, 8).Value = "Berries"
$ws.Cells.Item(23, 9).Value = 100101001
$ws.Cells.Item(23, 10).Value = "Arándano (blue)"
$ws.Cells.Item(23, 11).Value = "Sin especificar"
$ws.Cells.Item(23, 12).Value = "Primera"
$ws.Cells.Item(23, 13).Value = 180
$ws.Cells.Item(23, 14).Value = 4000
$ws.Cells.Item(23, 15).Value = 4000
$ws.Cells.Item(23, 16).Value = 4000
$ws.Cells.Item(23, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(23, 18).Value = "Provincia de Linares"
$ws.Cells.Item(23, 19).Value = 2000
$ws.Cells.Item(23, 20).Value = 2
